$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# Column C holds "Country WHS" codes stored as text (leading zeros matter,
# e.g. "07"), so force Text format before writing, then restore the
# original numeric display format ("0") afterwards without disturbing the
# already-written text values or the cell style index.
$ws.Range("C2:C5").NumberFormat = "@"

$ws.Range("C2").Value = "07"
$ws.Range("C3").Value = "43"
$ws.Range("C4").Value = "82"
$ws.Range("C5").Value = "82"

$ws.Range("C2:C5").NumberFormat = "0"
